# "Generate Report for Handback"
# Updates the Overview / zh-cn / de-de sheets of the localization-status
# report to reflect that the zh-cn and de-de handoffs have now been
# handed back (in sync with en-US).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both rows
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# Widen the (now longer) status columns
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
# Status column gets wider to fit the new text
$zhcn.Columns.Item(3).ColumnWidth = 29.15
# Latest Target File / Latest Handback File columns widen to 40
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# Status text (shared across the report)
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

# Row 2 (3028561f...) - fill in the Latest Target File / Latest Handback
# File / Latest Handback DateTime now that the handback happened
$zhcn.Range("J2").Value = "3028561f-dd41-4622-acdd-702562b5e347.e1a02ae33c1011caae370e836a078c443c18b0cc.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-07 02:30:42"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdd1866c147eea8d11b8f505b6b8d994bd634dd7/e2e/3028561f-dd41-4622-acdd-702562b5e347.md", "", "", "3028561f-dd41-4622-acdd-702562b5e347.md")

# Row 3 (ae869d6d...)
$zhcn.Range("J3").Value = "ae869d6d-2931-444f-b2bc-c8438c498020.ff524c2738ec79959d12817a873590bf694dba50.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-07 02:30:42"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdd1866c147eea8d11b8f505b6b8d994bd634dd7/e2e/ae869d6d-2931-444f-b2bc-c8438c498020.md", "", "", "ae869d6d-2931-444f-b2bc-c8438c498020.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# Row 2 (3028561f...)
$dede.Range("J2").Value = "3028561f-dd41-4622-acdd-702562b5e347.e1a02ae33c1011caae370e836a078c443c18b0cc.de-de.xlf"
$dede.Range("K2").Value = "2016-09-07 02:30:50"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdd1866c147eea8d11b8f505b6b8d994bd634dd7/e2e/3028561f-dd41-4622-acdd-702562b5e347.md", "", "", "3028561f-dd41-4622-acdd-702562b5e347.md")

# Row 3 (ae869d6d...)
$dede.Range("J3").Value = "ae869d6d-2931-444f-b2bc-c8438c498020.ff524c2738ec79959d12817a873590bf694dba50.de-de.xlf"
$dede.Range("K3").Value = "2016-09-07 02:30:50"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdd1866c147eea8d11b8f505b6b8d994bd634dd7/e2e/ae869d6d-2931-444f-b2bc-c8438c498020.md", "", "", "ae869d6d-2931-444f-b2bc-c8438c498020.md")
